$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "update automatically" date field from 17/02/2025 to
#    21/02/2025 everywhere it is cached: the Slide Master, every Custom
#    Layout (slide layout) and the Notes Master.
# ---------------------------------------------------------------------------

$oldDate = "17/02/2025"
$newDate = "21/02/2025"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout that hangs off the master
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# Notes Master (its date placeholder text only updates through the
# HeadersFooters facade in this environment)
$notesMaster = $p.NotesMaster
$nmDate = $notesMaster.HeadersFooters.DateAndTime
if ($nmDate.Text -eq $oldDate -or $nmDate.Text -eq "") {
    $nmDate.Text = $newDate
}

# ---------------------------------------------------------------------------
# 2) Slide 29 ("UI and Visualisation"): nudge "Picture 20" to the right
#    (95778 EMU -> 332623 EMU == 7.541575pt -> 26.19079pt).
# ---------------------------------------------------------------------------

$slide29 = $p.Slides.Item(29)
$pic = $slide29.Shapes.Item(2)
if ($pic.Name -eq "Picture 20") {
    $pic.Left = 332623 / 914400 * 72
}
